$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6
$ws.Cells.Item(6, 1).Value = 10034.51
$ws.Cells.Item(6, 2).Value = 9968.7199999999993
$ws.Cells.Item(6, 3).Value = 305.24
$ws.Cells.Item(6, 4).Value = 307.24
$ws.Cells.Item(6, 5).Value = $false
$ws.Cells.Item(6, 6).Value = 0.66
$ws.Cells.Item(6, 7).Value = 42613.766759259262
$ws.Cells.Item(6, 8).Value = $true

# Row 7
$ws.Cells.Item(7, 1).Value = 10052.57
$ws.Cells.Item(7, 2).Value = 10034.51
$ws.Cells.Item(7, 3).Value = 307.68
$ws.Cells.Item(7, 4).Value = 308.24
$ws.Cells.Item(7, 5).Value = $false
$ws.Cells.Item(7, 6).Value = 0.18
$ws.Cells.Item(7, 7).Value = 42614.674768518518
$ws.Cells.Item(7, 8).Value = $true

# Row 8
$ws.Cells.Item(8, 1).Value = 10048.549999999999
$ws.Cells.Item(8, 2).Value = 10052.57
$ws.Cells.Item(8, 3).Value = 307.95999999999998
$ws.Cells.Item(8, 4).Value = 307.83
$ws.Cells.Item(8, 5).Value = $false
$ws.Cells.Item(8, 6).Value = -0.04
$ws.Cells.Item(8, 7).Value = 42615.751851851855
$ws.Cells.Item(8, 8).Value = $false

# Ensure date/time style (style index 1, numFmtId 22) is applied to the new G cells,
# matching the existing column G formatting.
$ws.Range("G6:G8").NumberFormat = "m/d/yy h:mm"
